# Update the five "two-digit divided by one-digit" practice rows with a
# newly generated set of problems. Each populated row of the table is
# addressed by its (row, column) index and its run text is replaced
# in-place, which preserves run/paragraph formatting (font, size,
# justification) and avoids any ambiguity from duplicate values that
# appear as both old and new text elsewhere in the table.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row 1
$t.Cell(1,1).Range.Text = "45÷8="
$t.Cell(1,2).Range.Text = "14÷8="
$t.Cell(1,3).Range.Text = "19÷2="
$t.Cell(1,4).Range.Text = "59÷6="
$t.Cell(1,5).Range.Text = "75÷7="

# Row 5
$t.Cell(5,1).Range.Text = "79÷2="
$t.Cell(5,2).Range.Text = "76÷6="
$t.Cell(5,3).Range.Text = "22÷6="
$t.Cell(5,4).Range.Text = "99÷9="
$t.Cell(5,5).Range.Text = "65÷8="

# Row 9
$t.Cell(9,1).Range.Text = "47÷8="
$t.Cell(9,2).Range.Text = "18÷4="
$t.Cell(9,3).Range.Text = "74÷6="
$t.Cell(9,4).Range.Text = "21÷6="
$t.Cell(9,5).Range.Text = "26÷9="

# Row 13
$t.Cell(13,1).Range.Text = "62÷5="
$t.Cell(13,2).Range.Text = "42÷3="
$t.Cell(13,3).Range.Text = "14÷9="
$t.Cell(13,4).Range.Text = "17÷4="
$t.Cell(13,5).Range.Text = "70÷9="

# Row 17
$t.Cell(17,1).Range.Text = "18÷6="
$t.Cell(17,2).Range.Text = "93÷4="
$t.Cell(17,3).Range.Text = "60÷3="
$t.Cell(17,4).Range.Text = "65÷9="
$t.Cell(17,5).Range.Text = "82÷4="
